$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is plain numeric-looking text (e.g. "329.55") need a
# NumberFormat="@" round-trip so Excel stores them as text instead of a number,
# matching the original inline-string cell type. The cells original Style
# object is captured and restored afterwards so no visible style/format change
# is left behind on the cell.

$ws.Range("D2").Value = "30.211.95"
$ws.Range("E2").Value = "  +5.48%  "
$ws.Range("D3").Value = "1.914.32"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  -0.49%  "
$sTmp = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.55"
$ws.Range("D5").Style = $sTmp
$ws.Range("E5").Value = "  +4.80%  "
$sTmp = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = $sTmp
$ws.Range("E6").Value = "  -0.42%  "
$sTmp = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5192"
$ws.Range("D7").Style = $sTmp
$ws.Range("E7").Value = "  +2.24%  "
$sTmp = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4066"
$ws.Range("D8").Style = $sTmp
$ws.Range("E8").Value = "  +3.70%  "
$sTmp = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08494"
$ws.Range("D9").Style = $sTmp
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("E10").Value = "  +1.80%  "
$sTmp = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.77"
$ws.Range("D11").Style = $sTmp
$ws.Range("E11").Value = "  +1.13%  "
$sTmp = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.34"
$ws.Range("D12").Style = $sTmp
$ws.Range("E12").Value = "  +14.65%  "
$sTmp = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.452"
$ws.Range("D13").Style = $sTmp
$ws.Range("E13").Value = "  +4.16%  "
$ws.Range("D14").Value = "1.909.48"
$ws.Range("E14").Value = "  +1.73%  "
$sTmp = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.390"
$ws.Range("D15").Style = $sTmp
$ws.Range("E15").Value = "  +1.75%  "
$sTmp = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = $sTmp
$ws.Range("E16").Value = "  -0.55%  "
$sTmp = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.06"
$ws.Range("D17").Style = $sTmp
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("E18").Value = "  +1.20%  "
$sTmp = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06701"
$ws.Range("D19").Style = $sTmp
$ws.Range("E19").Value = "  -0.20%  "
$sTmp = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.52"
$ws.Range("D20").Style = $sTmp
$ws.Range("E20").Value = "  +4.98%  "
$ws.Range("E21").Value = "  -0.45%  "
$sTmp = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.018"
$ws.Range("D22").Style = $sTmp
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("D23").Value = "30.234.42"
$ws.Range("E23").Value = "  +5.46%  "
$sTmp = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.35"
$ws.Range("D24").Style = $sTmp
$ws.Range("E24").Value = "  +2.29%  "
$sTmp = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.232"
$ws.Range("D25").Style = $sTmp
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").Value = "2.136.90"
$ws.Range("E26").Value = "  +2.50%  "
$sTmp = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.44"
$ws.Range("D27").Style = $sTmp
$ws.Range("E27").Value = "  +3.65%  "
$sTmp = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.54"
$ws.Range("D28").Style = $sTmp
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("E29").Value = "  -0.60%  "
$sTmp = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.03"
$ws.Range("D30").Style = $sTmp
$ws.Range("E30").Value = "  +1.95%  "
$sTmp = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.104"
$ws.Range("D31").Style = $sTmp
$ws.Range("E31").Value = "  +5.57%  "
$sTmp = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1067"
$ws.Range("D32").Style = $sTmp
$ws.Range("E32").Value = "  +2.89%  "
$sTmp = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.018"
$ws.Range("D33").Style = $sTmp
$sTmp = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.648"
$ws.Range("D34").Style = $sTmp
$ws.Range("E34").Value = "  +0.69%  "
$sTmp = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02493"
$ws.Range("D35").Style = $sTmp
$ws.Range("E35").Value = "  +1.51%  "
$sTmp = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06577"
$ws.Range("D36").Style = $sTmp
$ws.Range("E36").Value = "  +0.47%  "
$sTmp = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2211"
$ws.Range("D37").Style = $sTmp
$ws.Range("E37").Value = "  +2.23%  "
$sTmp = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.195"
$ws.Range("D38").Style = $sTmp
$ws.Range("E38").Value = "  +3.05%  "
$sTmp = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.229"
$ws.Range("D39").Style = $sTmp
$ws.Range("E39").Value = "  +3.30%  "
$sTmp = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.91"
$ws.Range("D40").Style = $sTmp
$ws.Range("E40").Value = "  +6.96%  "
$sTmp = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.819"
$ws.Range("D41").Style = $sTmp
$sTmp = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6525"
$ws.Range("D42").Style = $sTmp
$ws.Range("E42").Value = "  +2.22%  "
$sTmp = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.237"
$ws.Range("D43").Style = $sTmp
$ws.Range("E43").Value = "  -0.22%  "
$sTmp = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6146"
$ws.Range("D44").Style = $sTmp
$ws.Range("E44").Value = "  +2.56%  "
$sTmp = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.34"
$ws.Range("D45").Style = $sTmp
$ws.Range("E45").Value = "  +2.25%  "
$sTmp = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.741"
$ws.Range("D46").Style = $sTmp
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("E47").Value = "  +3.53%  "
$ws.Range("E48").Value = "  +2.08%  "
$sTmp = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.54"
$ws.Range("D49").Style = $sTmp
$ws.Range("E49").Value = "  +1.86%  "
$sTmp = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.162"
$ws.Range("D50").Style = $sTmp
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("E51").Value = "  +4.37%  "
